$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1.4618836743572776
$ws.Range("B2").Value = -212.93231220423337
$ws.Range("C2").Value = -214.88162964715997
$ws.Range("D2").Value = -17.288273718396137
$ws.Range("E2").Value = -13.88282219047242
$ws.Range("F2").Value = 0.07950282610699369
$ws.Range("G2").Value = 1.5600414201183295
$ws.Range("A3").Value = -0.24383995651172938
$ws.Range("B3").Value = 33.6921607528605
$ws.Range("C3").Value = 33.795504650290674
$ws.Range("D3").Value = 2.7114605454817546
$ws.Range("E3").Value = 2.0303215827462866
$ws.Range("F3").Value = -0.004766771477098961
$ws.Range("G3").Value = -0.09334023668638827
$ws.Range("A4").Value = 0.0005971113936486139
$ws.Range("B4").Value = 22.75929259845534
$ws.Range("C4").Value = 22.84408668347703
$ws.Range("D4").Value = 2.75586954812938
$ws.Range("E4").Value = 2.7488713323011185
$ws.Range("F4").Value = -0.00010217023824321713
$ws.Range("G4").Value = 0.0004262573964502386
$ws.Range("A5").Value = -0.0037387724426624497
$ws.Range("B5").Value = 2.4694649505997033
$ws.Range("C5").Value = 3.4778713394214735
$ws.Range("D5").Value = 0.20038219812795371
$ws.Range("E5").Value = 0.22082235456114238
$ws.Range("F5").Value = -0.0008739415300380018
$ws.Range("G5").Value = -0.017021449704141774
$ws.Range("A6").Value = -0.000090278998200823965046384056
$ws.Range("B6").Value = -3.788269576945884
$ws.Range("C6").Value = -3.796155726562015
$ws.Range("D6").Value = -0.4588645547393215
$ws.Range("E6").Value = -0.4570368054086007
$ws.Range("F6").Value = 0.000011932263087399393522436419
$ws.Range("G6").Value = -0.000071042899408373111217254725
$ws.Range("A7").Value = 0.0020243537775510483
$ws.Range("B7").Value = -0.373711001567797
$ws.Range("C7").Value = -0.3698780535465467
$ws.Range("D7").Value = -0.02982068781611651
$ws.Range("E7").Value = -0.015740042053989135
$ws.Range("F7").Value = 0.00014288134640670714
$ws.Range("G7").Value = 0.002836908284023629
$ws.Range("A8").Value = -0.000006892149065428898568311158
$ws.Range("B8").Value = 0.09752824796314584
$ws.Range("C8").Value = 0.09695373591831853
$ws.Range("D8").Value = -0.006115181051572671
$ws.Range("E8").Value = -0.00594808282500277
$ws.Range("F8").Value = 0.00000124555599650084086575403
$ws.Range("G8").Value = 0.00000965236686389564488611735
$ws.Range("A9").Value = 0.000000390314176255957011678051
$ws.Range("B9").Value = -0.01634436779456557
$ws.Range("C9").Value = -0.01641152450306514
$ws.Range("D9").Value = 0.0010110503527073745
$ws.Range("E9").Value = 0.0009644603391038009
$ws.Range("F9").Value = -0.000000137058342703417536934582
$ws.Range("G9").Value = -0.000001608727810649274077100146